$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings; some are plain numeric-looking tokens
# (e.g. "252.25") that Excel would silently reinterpret as a Number if
# assigned directly. Prefixing with a literal leading apostrophe forces
# Excel to keep them as Text, same as typing them in manually.
$forceText = "'"

$ws.Range("D2").Value = '30.339.33'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '1.943.34'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = $forceText + '252.25'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = $forceText + '0.7255'
$ws.Range("E6").Value = '  -8.45%  '
$ws.Range("D7").Value = $forceText + '1.0000'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = $forceText + '0.3348'
$ws.Range("E8").Value = '  -4.54%  '
$ws.Range("D9").Value = $forceText + '28.92'
$ws.Range("E9").Value = '  +3.17%  '
$ws.Range("D10").Value = $forceText + '0.07440'
$ws.Range("E10").Value = '  +6.39%  '
$ws.Range("D11").Value = $forceText + '0.8219'
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").Value = $forceText + '0.08136'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").Value = '1.942.92'
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("D14").Value = $forceText + '5.503'
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("D15").Value = $forceText + '95.47'
$ws.Range("E15").Value = '  -4.81%  '
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").Value = $forceText + '0.000008425'
$ws.Range("E17").Value = '  +6.68%  '
$ws.Range("D18").Value = '30.364.70'
$ws.Range("E18").Value = '  -2.55%  '
$ws.Range("D19").Value = $forceText + '253.78'
$ws.Range("E19").Value = '  -7.08%  '
$ws.Range("D20").Value = $forceText + '5.907'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("D21").Value = '2.196.32'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("D22").Value = $forceText + '0.9999'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = $forceText + '1.001'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = $forceText + '7.002'
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").Value = $forceText + '9.893'
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("D26").Value = $forceText + '162.37'
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").Value = $forceText + '2.429'
$ws.Range("E27").Value = '  +4.92%  '
$ws.Range("D28").Value = $forceText + '19.42'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("D29").Value = $forceText + '0.1323'
$ws.Range("E29").Value = '  -12.41%  '
$ws.Range("D30").Value = $forceText + '1.574'
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("D31").Value = $forceText + '1.345'
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("D32").Value = $forceText + '4.471'
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("D33").Value = $forceText + '4.271'
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("D34").Value = $forceText + '0.05300'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").Value = $forceText + '1.314'
$ws.Range("E35").Value = '  +7.12%  '
$ws.Range("D36").Value = $forceText + '0.7615'
$ws.Range("E36").Value = '  -2.13%  '
$ws.Range("D37").Value = $forceText + '2.755'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = $forceText + '0.01997'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = $forceText + '2.852'
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").Value = $forceText + '81.49'
$ws.Range("D41").Value = $forceText + '6.625'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").Value = $forceText + '0.4591'
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("D43").Value = $forceText + '2.051'
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("D44").Value = $forceText + '0.8471'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = $forceText + '103.09'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("D47").Value = $forceText + '9.892'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").Value = $forceText + '7.519'
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("D49").Value = $forceText + '37.07'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = $forceText + '0.4217'
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = $forceText + '1.518'
$ws.Range("E51").Value = '  -0.39%  '
